# Weekly price report update: insert a new record at row 83 (shifting the
# existing rows 83-222 down to 84-223) for "Feria Lagunitas de Puerto Montt -
# Betarraga". This mirrors the commit's "semanal" (weekly) roll described in
# the commit message, where a newer data point is added at the top of the
# series and the sheet's used range grows by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 83:222 down to 84:223, leaving a blank row 83 behind.
$ws.Rows(83).Insert()

# Populate the new row 83 with the new weekly record.
$ws.Range("A83").Value = 4
$ws.Range("B83").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C83").Value = "Los Lagos"
$ws.Range("D83").Value = 44540
$ws.Range("E83").Value = 10
$ws.Range("F83").Value = 100114014
$ws.Range("G83").Value = "Betarraga"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 900
$ws.Range("L83").Value = 1000
$ws.Range("M83").Value = 950
$ws.Range("N83").Value = "`$/paquete 5 unidades"
$ws.Range("O83").Value = "Región del Maule"
$ws.Range("P83").Value = 190
$ws.Range("Q83").Value = 5
$ws.Range("R83").Value = "Hortaliza"
